$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Second matrix (remote pair programming), rows 14-20, cols B-H ---

# Row 14 (Felix): clear the "3rd Observer Ermal" label in B14 (self pairing cell),
# and update Felix/Magdalena minutes.
$ws.Range("B14").ClearContents()
$ws.Range("F14").Value = 240

# Row 15 (Ermal): mirror Felix/Ermal value from C14.
$ws.Range("B15").Formula = "=C14"

# Row 16 (Michael): mirror values from column D, update Michael/Jakob minutes.
$ws.Range("B16").Formula = "=D14"
$ws.Range("C16").Formula = "=D15"
$ws.Range("E16").Value = 360

# Row 17 (Jakob): mirror values from column E.
$ws.Range("B17").Formula = "=E14"
$ws.Range("D17").Formula = "=E16"

# Row 18 (Magdalena): mirror values from column F, clear observer label in F18.
$ws.Range("B18").Formula = "=F14"
$ws.Range("C18").Formula = "=F15"
$ws.Range("E18").Formula = "=E17"
$ws.Range("F18").ClearContents()

# Row 19 (Thomas): update Thomas/Florian minutes.
$ws.Range("H19").Value = 240

# Row 20 (Florian): mirror values from column H, update Thomas/Florian minutes.
$ws.Range("B20").Formula = "=H14"
$ws.Range("C20").Formula = "=H15"
$ws.Range("D20").Formula = "=H16"
$ws.Range("E20").Formula = "=H17"
$ws.Range("F20").Formula = "=H18"
$ws.Range("G20").Value = 240

# Leave the selection where the author last worked.
$ws.Range("J16").Select() | Out-Null

$wb.Save()
